$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "66.302.74"
$cell.ClearFormats()
$ws.Range("E2").Value = "  -2.70%  "
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "2.412.77"
$cell.ClearFormats()
$ws.Range("E3").Value = "  -4.79%  "
$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell.ClearFormats()
$ws.Range("E4").Value = "  +0.01%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "570.34"
$cell.ClearFormats()
$ws.Range("E5").Value = "  -4.00%  "
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "163.73"
$cell.ClearFormats()
$ws.Range("E6").Value = "  -7.60%  "
$ws.Range("E7").Value = "  +0.15%  "
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.499"
$cell.ClearFormats()
$ws.Range("E8").Value = "  -6.15%  "
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "2.411.64"
$cell.ClearFormats()
$ws.Range("E9").Value = "  -4.80%  "
$ws.Range("E10").Value = "  -8.64%  "
$ws.Range("E11").Value = "  -1.35%  "
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.322"
$cell.ClearFormats()
$ws.Range("E12").Value = "  -7.02%  "
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "4.73"
$cell.ClearFormats()
$ws.Range("E13").Value = "  -7.82%  "
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "24.55"
$cell.ClearFormats()
$ws.Range("E14").Value = "  -8.54%  "
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "65.975.29"
$cell.ClearFormats()
$ws.Range("E15").Value = "  -3.24%  "
$ws.Range("E16").Value = "  -8.37%  "
$ws.Range("E17").Value = "  -8.83%  "
$ws.Range("B28").Value = "WrappedeETH"
$ws.Range("C28").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "2.520.57"
$cell.ClearFormats()
$ws.Range("E28").Value = "  -5.08%  "
$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "0.0₃0868"
$cell.ClearFormats()
$ws.Range("E29").Value = "  -13.02%  "
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "7.62"
$cell.ClearFormats()
$ws.Range("E30").Value = "  -8.01%  "
$ws.Range("B31").Value = "Bittensor"
$ws.Range("C31").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "484.01"
$cell.ClearFormats()
$ws.Range("E31").Value = "  -10.70%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "1.74"
$cell.ClearFormats()
$ws.Range("E32").Value = "  -7.31%  "
$ws.Range("B33").Value = "FirstDigitalUSD"
$ws.Range("C33").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "0.999"
$cell.ClearFormats()
$ws.Range("E33").Value = "  -0.03%  "
$ws.Range("B34").Value = "Fetch.AI"
$ws.Range("C34").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "1.19"
$cell.ClearFormats()
$ws.Range("E34").Value = "  -10.93%  "
$ws.Range("B35").Value = "Monero"
$ws.Range("C35").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "156.13"
$cell.ClearFormats()
$ws.Range("E35").Value = "  -0.71%  "
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "0.114"
$cell.ClearFormats()
$ws.Range("E36").Value = "  -11.83%  "
$ws.Range("B37").Value = "WhiteBITCoin"
$ws.Range("C37").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "18.53"
$cell.ClearFormats()
$ws.Range("E37").Value = "  -0.94%  "
$ws.Range("B38").Value = "EthereumClassic"
$ws.Range("C38").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "18.07"
$cell.ClearFormats()
$ws.Range("E38").Value = "  -4.16%  "
$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "1.31"
$cell.ClearFormats()
$ws.Range("E39").Value = "  -10.16%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "1.63"
$cell.ClearFormats()
$ws.Range("E40").Value = "  -10.09%  "
$ws.Range("B41").Value = "PolygonEcosystemToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "0.319"
$cell.ClearFormats()
$ws.Range("E41").Value = "  -10.39%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "4.56"
$cell.ClearFormats()
$ws.Range("E42").Value = "  -12.41%  "
$ws.Range("B43").Value = "OKB"
$ws.Range("C43").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "39.06"
$cell.ClearFormats()
$ws.Range("E43").Value = "  -2.21%  "
$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "2.29"
$cell.ClearFormats()
$ws.Range("E44").Value = "  -10.35%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "136.86"
$cell.ClearFormats()
$ws.Range("E45").Value = "  -7.16%  "
$ws.Range("B46").Value = "Filecoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "3.40"
$cell.ClearFormats()
$ws.Range("E46").Value = "  -8.93%  "
$ws.Range("B47").Value = "ARBITRUM"
$ws.Range("C47").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "0.502"
$cell.ClearFormats()
$ws.Range("E47").Value = "  -10.57%  "
$ws.Range("B48").Value = "Optimism"
$ws.Range("C48").Value = "https://coinranking.com/coin/n1p-s_gm1+optimism-op"
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "1.55"
$cell.ClearFormats()
$ws.Range("E48").Value = "  -9.33%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "0.0717"
$cell.ClearFormats()
$ws.Range("E49").Value = "  -5.30%  "
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "0.572"
$cell.ClearFormats()
$ws.Range("E50").Value = "  -4.28%  "
$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "0.0₆0240"
$cell.ClearFormats()
$ws.Range("E51").Value = "  -13.52%  "
